$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 364, shifting existing rows 364-371 down to 365-372
$ws.Rows(364).Insert()

# Populate the newly inserted row 364 with the new data record
$ws.Range("A364").Value = 8
$ws.Range("B364").Value = "Terminal La Palmera de La Serena"
$ws.Range("C364").Value = "Coquimbo"
$ws.Range("D364").Value = 44628
$ws.Range("E364").Value = 4
$ws.Range("F364").Value = 100114001
$ws.Range("G364").Value = "Papa"
$ws.Range("H364").Value = "Asterix"
$ws.Range("I364").Value = "1a (cosecha)"
$ws.Range("J364").Value = 2540
$ws.Range("K364").Value = 9000
$ws.Range("L364").Value = 10000
$ws.Range("M364").Value = 9500
$ws.Range("N364").Value = "$/saco 25 kilos"
$ws.Range("O364").Value = "Región de Los Lagos"
$ws.Range("P364").Value = 380
$ws.Range("Q364").Value = 25
$ws.Range("R364").Value = "Hortaliza"
